{"js": "// Edit: expand the paper title from\n//   \"Review on Mechanical Issues and Driver Solutions of Industrial Servo Systems\"\n// to\n//   \"Review on Mechanical Issues and Driver Solutions of Industrial PMAC Servo\n//    Systems: Parameter Estimation and Auto-tuning Concepts\"\n// The canonical OOXML diff splits the (unchanged-formatting) title run into\n// three runs - this mirrors a user placing the caret mid-word (\"...Driver\n// Solut|ions...\") and typing the new wording, then appending \" Concepts\" at\n// the very end. We reproduce that exact run layout with Range.insertOoxml so\n// the first (untouched head) run keeps its original rsid attribute while the\n// two freshly authored runs come back bare, exactly like the source diff.\n\nconst ORIGINAL_TITLE =\n  \"Review on Mechanical Issues and Driver Solutions of Industrial Servo Systems\";\nconst HEAD_TEXT = \"Review on Mechanical Issues and Driver Solut\";\nconst MID_TEXT =\n  \"ions of Industrial PMAC Servo Systems: Parameter Estimation and Auto-tuning\";\nconst TAIL_TEXT = \" Concepts\";\n\nconst body = context.document.body;\nconst titleMatches = body.search(ORIGINAL_TITLE, { matchCase: true });\ntitleMatches.load(\"items\");\nawait context.sync();\n\nif (titleMatches.items.length === 0) {\n  throw new Error(\"Could not find the title text to update.\");\n}\n\nconst titleRange = titleMatches.items[0];\n\n// Run formatting shared by the original run and the two newly authored runs\n// (matches the <w:rPr> already present on the title run in before.docx).\nconst rFonts = `<w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/>`;\nconst rPr = `<w:rPr>${rFonts}<w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-GB\"/></w:rPr>`;\n\n// Flat-OPC OOXML package understood by Range.insertOoxml(). The first run\n// keeps the same rsid the original author run had; the two new runs (the\n// mid rewritten segment + the appended \" Concepts\") are plain <w:r> just\n// like freshly typed text.\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidRPr=\"00026D7C\">${rPr}<w:t>${HEAD_TEXT}</w:t></w:r>\n            <w:r>${rPr}<w:t>${MID_TEXT}</w:t></w:r>\n            <w:r>${rPr}<w:t xml:space=\"preserve\">${TAIL_TEXT}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntitleRange.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Edit: expand the paper title from\n#   \"Review on Mechanical Issues and Driver Solutions of Industrial Servo Systems\"\n# to\n#   \"Review on Mechanical Issues and Driver Solutions of Industrial PMAC Servo\n#    Systems: Parameter Estimation and Auto-tuning Concepts\"\n# The canonical OOXML diff splits the (unchanged-formatting) title run into\n# three runs - this mirrors a user placing the caret mid-word (\"...Driver\n# Solut|ions...\") and typing the new wording, then appending \" Concepts\" at\n# the very end. We reproduce that exact run layout via Range.InsertXML so the\n# first (untouched head) run keeps its original rsid attribute while the two\n# freshly authored runs come back bare, exactly like the source diff.\n\n$d = $word.ActiveDocument\n\n$originalTitle = \"Review on Mechanical Issues and Driver Solutions of Industrial Servo Systems\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($originalTitle)\nif (-not $found) {\n    throw \"Could not find the title text to update.\"\n}\n\n# Run formatting shared by the original run and the two newly authored runs.\n$rPr = '<w:rPr><w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\" w:cs=\"Times New Roman\"/><w:b/><w:sz w:val=\"28\"/><w:szCs w:val=\"28\"/><w:lang w:val=\"en-GB\"/></w:rPr>'\n\n# Flat-OPC OOXML package consumed by Range.InsertXML(). The first run keeps\n# the original author's rsid; the two new runs (rewritten mid segment plus\n# the appended \" Concepts\") come back as plain <w:r>, just like freshly\n# typed text.\n$flatOpc = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r w:rsidRPr=\"00026D7C\">$rPr<w:t>Review on Mechanical Issues and Driver Solut</w:t></w:r>\n            <w:r>$rPr<w:t>ions of Industrial PMAC Servo Systems: Parameter Estimation and Auto-tuning</w:t></w:r>\n            <w:r>$rPr<w:t xml:space=\"preserve\"> Concepts</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$rng.InsertXML($flatOpc)\n"}
